# ADD: Manipulation check, Promotion-based final results, exploration post-hoc.
#
# The regression-output table in columns B:E (estimate, std.error, statistic,
# p.value) was re-rounded from 3 decimal places down to 2 decimal places.
# Only the cells whose displayed value actually changes under the new
# rounding are touched here; cells already at (or unaffected by) 2-decimal
# rounding are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = 1.75;  "C2"  = 0.82; "D2"  = 2.12;  "E2"  = 0.03
    "B3"  = -0.01; "C3"  = 0.01
    "B4"  = -0.49; "C4"  = 0.22; "D4"  = -2.27; "E4"  = 0.02
    "C5"  = 0.07;  "D5"  = 0.45; "E5"  = 0.65
    "B6"  = -0.11; "C6"  = 0.18; "D6"  = -0.63; "E6"  = 0.53
    "B7"  = 0.01;  "C7"  = 0.09; "D7"  = 0.13;  "E7"  = 0.89
    "B8"  = -0.06; "C8"  = 0.09; "D8"  = -0.66; "E8"  = 0.51
    "B9"  = -0.03; "C9"  = 0.21; "D9"  = -0.14; "E9"  = 0.89
    "B10" = -0.29; "C10" = 0.21; "D10" = -1.33; "E10" = 0.18
    "B11" = 0.14;  "C11" = 0.24; "D11" = 0.57;  "E11" = 0.57
    "B12" = -0.19; "C12" = 0.3;  "D12" = -0.64; "E12" = 0.52
    "C13" = 0.3;   "D13" = -1.06; "E13" = 0.29
    "B14" = -0.17; "C14" = 0.3;  "D14" = -0.57; "E14" = 0.57
    "B15" = 0.15;  "C15" = 0.3;  "E15" = 0.63
    "B16" = 0.14;  "C16" = 0.31; "D16" = 0.46;  "E16" = 0.65
    "B17" = 0.08;  "C17" = 0.26; "D17" = 0.32;  "E17" = 0.75
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
